# "Generate Report for Archive" — refresh the localization-status report:
#   * the e2e sample file has moved from hand-off into active translation,
#     so its Status cells (shown on the Overview roll-up and on each
#     per-locale detail sheet) change from "Ready for handoff" to
#     "In Translation"
#   * the narrower status text lets the Status column(s) shrink to fit

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview roll-up: zh-cn / de-de status columns (E, F) for the one tracked file
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-locale detail sheets: Status column (C) for the same file
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Shrink the Status column(s) now that the text is shorter
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
